# case_deactivate.xlsx -- "Fix some issues with forms"
#
# - Bump the form_title / form_id on the "settings" sheet to the new
#   2023-02-09 revision.
# - Make "settings" the active/selected sheet (it was "choices" before).
# - Move the lingering selection on "survey" from C8 to A8.
# - Re-size column A on "settings" so the longer form_id/title values fit.

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$choices  = $wb.Worksheets.Item("choices")
$settings = $wb.Worksheets.Item("settings")

# --- Update the two strings that encode the form's title/id -------------
$settings.Range("A2").Value = "Household Exit Survey – 20230209"
$settings.Range("B2").Value = "case_deactivate_20230209"

# --- Widen column A on "settings" to fit the refreshed values -----------
$settings.Columns.Item(1).ColumnWidth = 29.5

# --- Leftover selection on "survey" moves from C8 to A8 ------------------
$null = $survey.Range("A8").Select()

# --- "settings" becomes the active tab (was "choices") -------------------
$null = $settings.Activate()
$null = $settings.Range("A2").Select()
